# Update the cryptos price/volume table with the latest scraped values.
# Price cells in column D are stored as text (not numbers) in the
# workbook, so values are entered with a leading apostrophe to force
# text entry, then the style is reset to "Normal" so no stray
# text-number style (quote-prefix) gets attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").Value = '''243.56'
$ws.Range("D2").Style = "Normal"

# Row 3 - OKB
$ws.Range("D3").Value = '''23.06'
$ws.Range("D3").Style = "Normal"

# Row 4 - HuobiToken
$ws.Range("D4").Value = '''5.404'
$ws.Range("D4").Style = "Normal"

# Row 5 - Cronos
$ws.Range("D5").Value = '''0.05985'
$ws.Range("D5").Style = "Normal"

# Row 7 - KuCoinToken
$ws.Range("D7").Value = '''6.502'
$ws.Range("D7").Style = "Normal"

# Row 8 - MXToken
$ws.Range("D8").Value = '''0.8125'
$ws.Range("D8").Style = "Normal"

# Row 9 - FTXToken
$ws.Range("D9").Value = '''0.9292'
$ws.Range("D9").Style = "Normal"

# Row 10 - WazirX
$ws.Range("D10").Value = '''0.1425'
$ws.Range("D10").Style = "Normal"

# Row 11 - MandalaExchangeToken
$ws.Range("D11").Value = '''0.07421'
$ws.Range("D11").Style = "Normal"

# Row 12 - LiechtensteinCryptoassetsExchange
$ws.Range("D12").Value = '''0.03317'
$ws.Range("D12").Style = "Normal"

# Row 13 - BitrueCoin
$ws.Range("D13").Value = '''0.03064'
$ws.Range("D13").Style = "Normal"

# Row 15 - MCDex
$ws.Range("D15").Value = '''3.849'
$ws.Range("D15").Style = "Normal"

# Row 16 - BitForexToken
$ws.Range("D16").Value = '''0.001582'
$ws.Range("D16").Style = "Normal"

# Row 17 - CoinExToken
$ws.Range("D17").Value = '''0.04704'
$ws.Range("D17").Style = "Normal"

# Row 18 - One
$ws.Range("D18").Value = '''0.0005900'
$ws.Range("D18").Style = "Normal"

# Row 19 - TigerCash
$ws.Range("D19").Value = '''0.005953'
$ws.Range("D19").Style = "Normal"

# Row 20 - BitKan
$ws.Range("D20").Value = '''0.001271'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '19BitKanKANBestin24h'

# Row 21 - HotbitToken
$ws.Range("D21").Value = '''0.004880'
$ws.Range("D21").Style = "Normal"

# Row 22 - now NitroEx (was UpBots)
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").Value = '''0.00007998'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '21NitroExNTX'

# Row 23 - now LEO (was NitroEx)
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").Value = '''3.571'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '22LEOLEO'

# Row 24 - now BTSEToken (was LEO)
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").Value = '''2.133'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '23BTSETokenBTSE'

# Row 25 - now BitpandaEcosystemToken (was BTSEToken)
$ws.Range("B25").Value = 'BitpandaEcosystemToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D25").Value = '''0.3240'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '24BitpandaEcosystemTokenBEST'

# Row 26 - now ProBitToken (was BitpandaEcosystemToken)
$ws.Range("B26").Value = 'ProBitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D26").Value = '''0.1331'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '25ProBitTokenPROB'

# Row 27 - now UpBots (was ProBitToken)
$ws.Range("B27").Value = 'UpBots'
$ws.Range("C27").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D27").Value = '''0.0002339'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '26UpBotsUBXT'

# Row 40 - IDEX
$ws.Range("D40").Value = '''0.03955'
$ws.Range("D40").Style = "Normal"

# Row 41 - now BKEXToken (was KickToken)
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = '''0.1077'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '40BKEXTokenBKK'

# Row 42 - CEJI
$ws.Range("D42").Value = '''0.002659'
$ws.Range("D42").Style = "Normal"

# Row 43 - now KickToken (was BKEXToken)
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = '''0.003071'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '42KickTokenKICKWorstin24h'

# Row 44 - LocalTraders
$ws.Range("D44").Value = '''0.009079'
$ws.Range("D44").Style = "Normal"

# Row 45 - CoinLion
$ws.Range("D45").Value = '''0.00005068'
$ws.Range("D45").Style = "Normal"

# Row 46 - Kangarootoken
$ws.Range("D46").Value = '''0.00000000750'
$ws.Range("D46").Style = "Normal"

# Row 47 - CoinbaseStockToken
$ws.Range("D47").Value = '''0.6700'
$ws.Range("D47").Style = "Normal"

# Row 48 - BOLO
$ws.Range("D48").Value = '''0.002266'
$ws.Range("D48").Style = "Normal"

# Row 49 - CryptobidCoin
$ws.Range("D49").Value = '''0.00002100'
$ws.Range("D49").Style = "Normal"

# Row 50 - SpecialPowerGold
$ws.Range("D50").Value = '''0.0002000'
$ws.Range("D50").Style = "Normal"
